# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (AC1) onto the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Data rows 2 through 49 get the same Wins/Losses/Ties values
$lastRow = 49
$ws.Range("AD2:AD$lastRow").Value = 68
$ws.Range("AE2:AE$lastRow").Value = 94
$ws.Range("AF2:AF$lastRow").Value = 0
